# Update the cryptos list worksheet with refreshed Price (col D) and
# Volume(1h) (col E) figures, and fix the Algorand / InternetComputer
# row ordering (rows 38 & 39 swapped, with refreshed values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where both the Price (D) and Volume(1h) (E) columns changed.
$rowUpdates = @{
    2  = @("23.415.53",  "  -1.35%  ")
    3  = @("1.645.63",   "  -0.49%  ")
    4  = @("1.004",      "  +0.48%  ")
    5  = @("1.002",      "  +0.25%  ")
    6  = @("299.37",     "  -1.67%  ")
    7  = @("0.3792",     "  -0.68%  ")
    8  = @("0.3557",     "  -1.58%  ")
    9  = @("50.05",      "  -3.07%  ")
    10 = @("0.08101",    "  -1.80%  ")
    11 = @("1.217",      "  -2.72%  ")
    12 = @("1.002",      "  +0.16%  ")
    13 = @("21.99",      "  -2.84%  ")
    14 = @("6.391",      "  -2.40%  ")
    15 = @("7.347",      "  -0.90%  ")
    16 = @("0.00001194", "  -3.30%  ")
    17 = @("1.651.89",   "  -0.48%  ")
    18 = @("97.42",      "  +0.32%  ")
    19 = @("0.06954",    "  -0.31%  ")
    20 = @("6.753",      "  -0.83%  ")
    21 = @("17.25",      "  -2.66%  ")
    23 = @("12.41",      "  -1.75%  ")
    24 = @("23.461.02",  "  -1.11%  ")
    25 = @("2.505",      "  -1.55%  ")
    26 = @("2.903",      "  -6.01%  ")
    27 = @("20.87",      "  -2.15%  ")
    28 = @("153.63",     "  +0.61%  ")
    29 = @("5.200",      "  -0.63%  ")
    30 = @("132.73",     "  -2.07%  ")
    31 = @("1.834.61",   "  -0.32%  ")
    32 = @("6.913",      "  +0.08%  ")
    33 = @("2.116",      "  +0.77%  ")
    34 = @("11.78",      "  -2.80%  ")
    35 = @("1.013",      "  -7.10%  ")
    36 = @("0.02719",    "  -3.58%  ")
    37 = @("0.08733",    "  -1.26%  ")
    40 = @("13.06",      "  +1.56%  ")
    41 = @("0.06776",    "  -4.03%  ")
    42 = @("0.6880",     "  -2.83%  ")
    43 = @("1.309",      "  -2.47%  ")
    46 = @("0.6378",     "  -2.33%  ")
    47 = @("2.253",      "  -4.03%  ")
    48 = @("3.921",      "  -1.60%  ")
    49 = @("0.07722",    "  -3.43%  ")
    50 = @("127.28",     "  -0.77%  ")
    51 = @("1.148",      "  -3.78%  ")
}

# Price values (column D) are stored as plain text in the source sheet
# (e.g. "23.478.28", "0.3792", "5.200") rather than numbers, so force the
# column to Text format before writing -- otherwise Excel auto-detects
# numeric-looking strings and silently converts/truncates them (e.g.
# "5.200" -> 5.2).
$ws.Range("D2:D51").NumberFormat = "@"

foreach ($r in $rowUpdates.Keys) {
    $vals = $rowUpdates[$r]
    $ws.Cells.Item($r, 4).Value = $vals[0]
    $ws.Cells.Item($r, 5).Value = $vals[1]
}

# Rows where only the Volume(1h) (E) column changed.
$eOnlyUpdates = @{
    44 = "  -2.87%  "
    45 = "  +0.17%  "
}

foreach ($r in $eOnlyUpdates.Keys) {
    $ws.Cells.Item($r, 5).Value = $eOnlyUpdates[$r]
}

# Rows 38 and 39: Algorand and InternetComputer(DFINITY) swap order,
# each also getting refreshed Price / Volume(1h) figures.
$ws.Cells.Item(38, 2).Value = "Algorand"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(38, 4).Value = "0.2430"
$ws.Cells.Item(38, 5).Value = "  -3.82%  "

$ws.Cells.Item(39, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(39, 4).Value = "5.928"
$ws.Cells.Item(39, 5).Value = "  -3.02%  "

# Restore the column's visual/style state (drop the temporary text
# format) now that every text-like price value has been written safely.
$ws.Range("D2:D51").Style = "Normal"
